# The deck's single red-outline "highlight" rectangles (Rectangle 24 / id 25
# and Rectangle 216 / id 217) that were overlaid on slide 1 are removed -
# mirrors the author deleting the two annotation boxes before re-saving.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$idsToDelete = @(25, 217)

for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $shape = $s.Shapes.Item($i)
    if ($idsToDelete -contains $shape.Id) {
        $shape.Delete()
    }
}
